$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "durata" column (header + values) in column E
$ws.Range("E1").Value = "durata"
$ws.Range("E2").Value = 120
$ws.Range("E3").Value = 150
$ws.Range("E4").Value = 138

# Update selection to mirror the authored change (E5, just below the data)
$ws.Range("E5").Select()
